$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 13, shifting the old rows 13..48 down to 14..49.
$ws.Rows.Item(13).Insert()

# The insert operation stamped an empty, styled A13 cell (inherited from the
# row above). Row 13 has no content in column A in the target layout, so
# drop it entirely.
$ws.Range("A13").Clear()

# Give the new B13/C13 cells the same formatting as their column neighbours
# (style 2 for column B, style 3 for column C) by copying format from the
# row below (which still carries the original styling).
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# --- Content updates ---

# Objetivos / Objectives body text (row 10)
$objetivos = "Fornecer oportunidade de realização de treinamento profissional de Engenharia Ambiental em empresa ou instituição sob supervisão de docente do Departamento de Ciências Básicas e Ambientais da EEL. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Docentes responsáveis (now row 13, no label cell, just the name)
$docente = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# Programa resumido / Short syllabus body text (row 14)
$resumido = "Processo seletivo. Plano de trabalho específico. Realização do estágio. Relatório final."
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# Programa / Syllabus body text (row 16)
$programa = "Participação do aluno em processo seletivo de empresas, instituições de pesquisa ou no setor acadêmico. O estágio realizado sob a supervisão de docente designado pelo Departamento de Ciências Básicas e Ambientais da Escola de Engenharia de Lorena. O conteúdo será estabelecido no Plano de Trabalho entre o supervisor responsável pelo Estágio e o docente supervisor. Apresentação de relatório final sobre as atividades desenvolvidas no estágio."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Método body text (row 19) - "Supervisão..." text, shifted up from its old row
$metodo = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Critério body text (row 20) - "A nota final..." text, shifted up
$criterio = "A nota final será baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Norma de recuperação body text (row 21) - "Devido às..." text, shifted up
$norma = "Devido às características da disciplina, não será oferecida recuperação."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Bibliografia body text (row 22) - new text "Não há."
$biblio = "Não há."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
